$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the claim number in row 2 (E2) - must stay text (leading apostrophe forces text entry)
$ws.Range("E2").Value = "'1120194100412"

# Insert two new rows before the old row 3, pushing it down to row 5.
# This also copies row 2's cell formatting (styles) down into the new rows.
$ws.Range("A3:A4").EntireRow.Insert()

# Fill new row 3 - same environment/url/user as row 2, new claim number
$ws.Range("A3").Value = "ssurgwsoadev4-oci.opc.oracleoutsourcing.com"
$ws.Range("B3").Value = "https://ssurgwsoadev4-oci.opc.oracleoutsourcing.com/cc/ClaimCenter.do"
$ws.Range("C3").Value = "dgariffo"
$ws.Range("D3").Value = "silverarrow"
$ws.Range("E3").Value = "'1220194200667"

# Fill new row 4 - same environment/url as row 2, new user, new claim number
$ws.Range("A4").Value = "ssurgwsoadev4-oci.opc.oracleoutsourcing.com"
$ws.Range("B4").Value = "https://ssurgwsoadev4-oci.opc.oracleoutsourcing.com/cc/ClaimCenter.do"
$ws.Range("C4").Value = "arubino"
$ws.Range("D4").Value = "silverarrow"
$ws.Range("E4").Value = "'0420194406717"

# The row insert kept the original hyperlink anchored at B3 instead of following the
# data down to B5 (old row 3). Recreate it in the right place.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B5"), "https://preproducciongestion.segurossura.com.ar/cc/ClaimCenter.do")
$ws.Range("B5").Style = "Hipervínculo"

# Update the active selection to match the edited rows
$ws.Range("A4:B4").Select()
